# Add thêm nhân sự Nguyễn Hữu Quang
# Updates the "Lương" (salary) worksheet: one more công (workday) added to
# CẦN THƠ, which ripples into the base-salary and total-salary figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B1").Value = 2
$ws.Range("B2").Value = 31
$ws.Range("B3").Value = 2214285.714285715

$ws.Range("B13").Value = 3321428.571428572
$ws.Range("B23").Value = 3321428.571428572

$ws.Range("B31").Value = 2214285.714285715
$ws.Range("B32").Value = 3321428.571428572
$ws.Range("B33").Value = 3321428.571428572
$ws.Range("B34").Value = 8857142.857142858
